# Auto update: 2025-12-05 19:04:26
# Updates the final score (K) and MACRO_SCORE (N) columns on rows 2-7
# to reflect refreshed analysis numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final score (column K / "최종점수")
$ws.Range("K2").Value = 62.5
$ws.Range("K3").Value = 54.3
$ws.Range("K4").Value = 50.5
$ws.Range("K5").Value = 47.7
$ws.Range("K6").Value = 39.7
$ws.Range("K7").Value = 39.7

# MACRO_SCORE (column N), same new value for all rows 2-7
$ws.Range("N2:N7").Value = 51.15965480231979
